$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("Add more client pages to use APIs (Add new User/Book)"): work has
# started locally (docker / docker-compose) -> record the actual start-date
# and flip status from "Not started" to "In-Progress". Pick up the green
# "In-Progress" look from G7 before G7's own status changes below.
$ws.Range("G7").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = "In-Progress"

$ws.Range("E8").Value2 = 43741
$ws.Range("E8").NumberFormat = $ws.Range("D8").NumberFormat

# Row 7 ("Deploy docker local"): task is now finished -> record the actual
# end-date and flip status to "Completed" (reuse the style already used by
# the other completed rows, e.g. G3).
$ws.Range("G3").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = "Completed"

$ws.Range("F7").Value2 = 43740
$ws.Range("F7").NumberFormat = $ws.Range("E7").NumberFormat

$excel.CutCopyMode = $false

# Reflect the final selected cell from the authored workbook.
$ws.Range("E9").Select()
